$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2-72 all get updated from serial 45203 to 45204
$ws.Range("C2:C72").Value = 45204
